$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ---- Sheet ALC ----
# Row 17
$ws_ALC.Range("H17").Value = 3166655.8
$ws_ALC.Range("J17").Value = 3166655.8
$ws_ALC.Range("L17").Value = 9499967.399999999
$ws_ALC.Range("N17").Value = -9500303.399999999

# Row 33
$ws_ALC.Range("H33").Value = 13892936
$ws_ALC.Range("I33").Value = 19231528
$ws_ALC.Range("K33").Value = 19231528
$ws_ALC.Range("M33").Value = -19231299

# Row 64
$ws_ALC.Range("H64").Value = 5999.909
$ws_ALC.Range("J64").Value = 6999.875
$ws_ALC.Range("L64").Value = 6999.875
$ws_ALC.Range("N64").Value = -7495.875

# Row 67
$ws_ALC.Range("H67").Value = 5999.909
$ws_ALC.Range("J67").Value = 6999.875
$ws_ALC.Range("L67").Value = 6999.875
$ws_ALC.Range("N67").Value = -8715.875

# Row 68
$ws_ALC.Range("H68").Value = 0
$ws_ALC.Range("J68").Value = 0
$ws_ALC.Range("N68").Value = 0
$ws_ALC.Range("L68").ClearContents()

# Row 71
$ws_ALC.Range("H71").Value = 0
$ws_ALC.Range("J71").Value = 0
$ws_ALC.Range("N71").Value = 0
$ws_ALC.Range("L71").ClearContents()

# Row 101
$ws_ALC.Range("H101").Value = 1114.6
$ws_ALC.Range("I101").Value = 1780.2
$ws_ALC.Range("K101").Value = 5340.6
$ws_ALC.Range("M101").Value = -3718.6

# Row 113
$ws_ALC.Range("I113").Value = 10970.5
$ws_ALC.Range("J113").Value = 5700
$ws_ALC.Range("K113").Value = 10970.5
$ws_ALC.Range("L113").Value = 5700
$ws_ALC.Range("M113").Value = -7716.5
$ws_ALC.Range("N113").Value = -12208

# Row 138
$ws_ALC.Range("H138").Value = 4506.0454
$ws_ALC.Range("I138").Value = 4373.4
$ws_ALC.Range("J138").Value = 4616.5835
$ws_ALC.Range("K138").Value = 13120.2
$ws_ALC.Range("L138").Value = 13849.7505
$ws_ALC.Range("M138").Value = -7980.199999999999
$ws_ALC.Range("N138").Value = -24129.7505

# ---- Sheet ARM ----
# Row 45
$ws_ARM.Range("H45").Value = 3466.7222
$ws_ARM.Range("I45").Value = 3305.1333
$ws_ARM.Range("K45").Value = 3305.1333
$ws_ARM.Range("M45").Value = -2928.1333

# Row 61
$ws_ARM.Range("H61").Value = 9664.208000000001
$ws_ARM.Range("I61").Value = 5167.5356
$ws_ARM.Range("J61").Value = 15959.55
$ws_ARM.Range("K61").Value = 5167.5356
$ws_ARM.Range("L61").Value = 15959.55
$ws_ARM.Range("M61").Value = -4955.5356
$ws_ARM.Range("N61").Value = -16383.55

# Row 122
$ws_ARM.Range("H122").Value = 2685.0344
$ws_ARM.Range("I122").Value = 1789.5897
$ws_ARM.Range("J122").Value = 4523.0527
$ws_ARM.Range("K122").Value = 5368.7691
$ws_ARM.Range("L122").Value = 13569.1581
$ws_ARM.Range("M122").Value = -2918.7691
$ws_ARM.Range("N122").Value = -18469.1581

# Row 136
$ws_ARM.Range("H136").Value = 9664.208000000001
$ws_ARM.Range("I136").Value = 5167.5356
$ws_ARM.Range("J136").Value = 15959.55
$ws_ARM.Range("K136").Value = 15502.6068
$ws_ARM.Range("L136").Value = 47878.64999999999
$ws_ARM.Range("M136").Value = -12952.6068
$ws_ARM.Range("N136").Value = -52978.64999999999

# ---- Sheet BSM ----
# Row 105
$ws_BSM.Range("H105").Value = 4575
$ws_BSM.Range("I105").Value = 4400
$ws_BSM.Range("K105").Value = 4400
$ws_BSM.Range("M105").Value = -2653

# Row 107
$ws_BSM.Range("H107").Value = 3837.5625
$ws_BSM.Range("I107").Value = 4206.846
$ws_BSM.Range("K107").Value = 4206.846
$ws_BSM.Range("M107").Value = -2286.846

# Row 134
$ws_BSM.Range("H134").Value = 11385.964
$ws_BSM.Range("I134").Value = 3515.611
$ws_BSM.Range("K134").Value = 10546.833
$ws_BSM.Range("M134").Value = -8011.832999999999

# ---- Sheet CRP ----
# Row 58
$ws_CRP.Range("H58").Value = 17909.656
$ws_CRP.Range("I58").Value = 7467.143
$ws_CRP.Range("J58").Value = 26031.611
$ws_CRP.Range("K58").Value = 7467.143
$ws_CRP.Range("L58").Value = 26031.611
$ws_CRP.Range("M58").Value = -7264.143
$ws_CRP.Range("N58").Value = -26437.611

# Row 70
$ws_CRP.Range("H70").Value = 16498.75
$ws_CRP.Range("I70").Value = 15995
$ws_CRP.Range("K70").Value = 15995
$ws_CRP.Range("M70").Value = -15680

# Row 73
$ws_CRP.Range("H73").Value = 16498.75
$ws_CRP.Range("I73").Value = 15995
$ws_CRP.Range("K73").Value = 15995
$ws_CRP.Range("M73").Value = -14903

# Row 94
$ws_CRP.Range("H94").Value = 1875.8334
$ws_CRP.Range("I94").Value = 2002.6666
$ws_CRP.Range("J94").Value = 1749
$ws_CRP.Range("K94").Value = 2002.6666
$ws_CRP.Range("L94").Value = 1749
$ws_CRP.Range("M94").Value = -1551.6666
$ws_CRP.Range("N94").Value = -2651

# Row 134
$ws_CRP.Range("H134").Value = 30310190
$ws_CRP.Range("I134").Value = 2741.8462
$ws_CRP.Range("J134").Value = 50010030
$ws_CRP.Range("K134").Value = 8225.5386
$ws_CRP.Range("L134").Value = 150030090
$ws_CRP.Range("M134").Value = -5690.5386
$ws_CRP.Range("N134").Value = -150035160

# Row 136
$ws_CRP.Range("H136").Value = 17909.656
$ws_CRP.Range("I136").Value = 7467.143
$ws_CRP.Range("J136").Value = 26031.611
$ws_CRP.Range("K136").Value = 22401.429
$ws_CRP.Range("L136").Value = 78094.833
$ws_CRP.Range("M136").Value = -19851.429
$ws_CRP.Range("N136").Value = -83194.833

# ---- Sheet CUL ----
# Row 4
$ws_CUL.Range("H4").Value = 27079934
$ws_CUL.Range("J4").Value = 111709.664
$ws_CUL.Range("L4").Value = 335128.992
$ws_CUL.Range("N4").Value = -335352.992

# Row 8
$ws_CUL.Range("H8").Value = 2778049.2
$ws_CUL.Range("I8").Value = 2778049.2
$ws_CUL.Range("K8").Value = 8334147.600000001
$ws_CUL.Range("M8").Value = -8334008.600000001

# Row 15
$ws_CUL.Range("H15").Value = 67.3
$ws_CUL.Range("I15").Value = 63.5
$ws_CUL.Range("J15").Value = 69.833336
$ws_CUL.Range("K15").Value = 190.5
$ws_CUL.Range("L15").Value = 209.500008
$ws_CUL.Range("M15").Value = -50.5
$ws_CUL.Range("N15").Value = -489.500008

# Row 98
$ws_CUL.Range("H98").Value = 12090.333
$ws_CUL.Range("J98").Value = 12090.333
$ws_CUL.Range("L98").Value = 36270.999
$ws_CUL.Range("N98").Value = -39266.999

# Row 114
$ws_CUL.Range("H114").Value = 2286.6667
$ws_CUL.Range("I114").Value = 1180.3334
$ws_CUL.Range("J114").Value = 4499.3335
$ws_CUL.Range("K114").Value = 3541.0002
$ws_CUL.Range("L114").Value = 13498.0005
$ws_CUL.Range("M114").Value = -287.0001999999999
$ws_CUL.Range("N114").Value = -20006.0005

# Row 122
$ws_CUL.Range("H122").Value = 13454728
$ws_CUL.Range("J122").Value = 3550355.2
$ws_CUL.Range("L122").Value = 31953196.8
$ws_CUL.Range("N122").Value = -31958096.8

# Row 131
$ws_CUL.Range("H131").Value = 1476.93
$ws_CUL.Range("J131").Value = 1483.6123
$ws_CUL.Range("L131").Value = 4450.8369
$ws_CUL.Range("N131").Value = -14530.8369

# Row 137
$ws_CUL.Range("H137").Value = 2089.611
$ws_CUL.Range("I137").Value = 1815.1666
$ws_CUL.Range("J137").Value = 2638.5
$ws_CUL.Range("K137").Value = 5445.4998
$ws_CUL.Range("L137").Value = 7915.5
$ws_CUL.Range("M137").Value = -345.4997999999996
$ws_CUL.Range("N137").Value = -18115.5

# Row 138
$ws_CUL.Range("H138").Value = 4149.8965
$ws_CUL.Range("I138").Value = 1305
$ws_CUL.Range("K138").Value = 3915
$ws_CUL.Range("M138").Value = 1225

# Row 139
$ws_CUL.Range("H139").Value = 4870.4546
$ws_CUL.Range("I139").Value = 4657.5
$ws_CUL.Range("K139").Value = 13972.5
$ws_CUL.Range("M139").Value = -8832.5

# Row 140
$ws_CUL.Range("H140").Value = 2530.5
$ws_CUL.Range("I140").Value = 1358.25
$ws_CUL.Range("K140").Value = 4074.75
$ws_CUL.Range("M140").Value = 1105.25

# Row 141
$ws_CUL.Range("H141").Value = 4680.727
$ws_CUL.Range("I141").Value = 1092.375
$ws_CUL.Range("K141").Value = 3277.125
$ws_CUL.Range("M141").Value = 1902.875

# ---- Sheet GSM ----
# Row 15
$ws_GSM.Range("H15").Value = 6217.093
$ws_GSM.Range("J15").Value = 6217.093
$ws_GSM.Range("L15").Value = 6217.093
$ws_GSM.Range("N15").Value = -6793.093

# Row 23
$ws_GSM.Range("H23").Value = 2000
$ws_GSM.Range("I23").Value = 0
$ws_GSM.Range("K23").Value = 0
$ws_GSM.Range("M23").ClearContents()

# Row 81
$ws_GSM.Range("H81").Value = 6217.093
$ws_GSM.Range("J81").Value = 6217.093
$ws_GSM.Range("L81").Value = 6217.093
$ws_GSM.Range("N81").Value = -8213.093000000001

# Row 84
$ws_GSM.Range("H84").Value = 6217.093
$ws_GSM.Range("J84").Value = 6217.093
$ws_GSM.Range("L84").Value = 18651.279
$ws_GSM.Range("N84").Value = -28635.279

# Row 113
$ws_GSM.Range("H113").Value = 65798.42999999999
$ws_GSM.Range("I113").Value = 91287.8
$ws_GSM.Range("J113").Value = 2075
$ws_GSM.Range("K113").Value = 91287.8
$ws_GSM.Range("L113").Value = 2075
$ws_GSM.Range("M113").Value = -89117.8
$ws_GSM.Range("N113").Value = -6415

# Row 122
$ws_GSM.Range("H122").Value = 3343.2778
$ws_GSM.Range("I122").Value = 1361.48
$ws_GSM.Range("K122").Value = 4084.44
$ws_GSM.Range("M122").Value = -1634.44

# ---- Sheet LTW ----
# Row 40
$ws_LTW.Range("H40").Value = 6078.4814
$ws_LTW.Range("I40").Value = 3831.3044
$ws_LTW.Range("K40").Value = 3831.3044
$ws_LTW.Range("M40").Value = -3695.3044

# Row 55
$ws_LTW.Range("H55").Value = 1802.619
$ws_LTW.Range("I55").Value = 658
$ws_LTW.Range("J55").Value = 2260.4666
$ws_LTW.Range("K55").Value = 658
$ws_LTW.Range("L55").Value = 2260.4666
$ws_LTW.Range("M55").Value = -485
$ws_LTW.Range("N55").Value = -2606.4666

# Row 122
$ws_LTW.Range("H122").Value = 6494.0464
$ws_LTW.Range("I122").Value = 5253.24
$ws_LTW.Range("K122").Value = 15759.72
$ws_LTW.Range("M122").Value = -13309.72

# Row 136
$ws_LTW.Range("H136").Value = 21148.459
$ws_LTW.Range("I136").Value = 33588.668
$ws_LTW.Range("J136").Value = 13684.333
$ws_LTW.Range("K136").Value = 100766.004
$ws_LTW.Range("L136").Value = 41052.999
$ws_LTW.Range("M136").Value = -98216.00399999999
$ws_LTW.Range("N136").Value = -46152.999

# ---- Sheet WVR ----
# Row 105
$ws_WVR.Range("H105").Value = 58481
$ws_WVR.Range("J105").Value = 58481
$ws_WVR.Range("L105").Value = 58481
$ws_WVR.Range("N105").Value = -65469

# Row 107
$ws_WVR.Range("H107").Value = 3409.476
$ws_WVR.Range("I107").Value = 1147.4
$ws_WVR.Range("J107").Value = 9064.666999999999
$ws_WVR.Range("K107").Value = 3442.2
$ws_WVR.Range("L107").Value = 27194.001
$ws_WVR.Range("M107").Value = -1522.2
$ws_WVR.Range("N107").Value = -31034.001

# Row 113
$ws_WVR.Range("H113").Value = 2516.125
$ws_WVR.Range("I113").Value = 2265.75
$ws_WVR.Range("K113").Value = 6797.25
$ws_WVR.Range("M113").Value = -4627.25

# Row 126
$ws_WVR.Range("H126").Value = 12401.421
$ws_WVR.Range("I126").Value = 14017.357
$ws_WVR.Range("K126").Value = 42052.071
$ws_WVR.Range("M126").Value = -39582.071

# Row 132
$ws_WVR.Range("H132").Value = 7924.1206
$ws_WVR.Range("I132").Value = 3158.3125
$ws_WVR.Range("K132").Value = 9474.9375
$ws_WVR.Range("M132").Value = -6944.9375
